$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: ADF4355-3BCPZ -------------------------------------------------
$ws.Range("A6").Value = "ADF4355-3BCPZ"
$ws.Range("B6").Value = "32-QFN (5x5)"
$ws.Range("C6").Value = "ADF4355-3BCPZ 32-QFN (5x5)"
$ws.Range("D6").Value = "ADF4355-3BCPZ"
$ws.Range("E6").Value = "ADF4355-3BCPZ"
$ws.Range("F6").Value = "32-QFN (5x5)"
$ws.Range("G6").Value = "Analog Devices"

# --- Row 7: ADF5356 --------------------------------------------------------
$ws.Range("A7").Value = "ADF5356"
$ws.Range("B7").Value = "32-QFN (5x5)"
$ws.Range("C7").Value = "ADF5356 32-QFN (5x5)"
$ws.Range("D7").Value = "ADF5356"
$ws.Range("E7").Value = "ADF5356"
$ws.Range("F7").Value = "32-QFN (5x5)"
$ws.Range("G7").Value = "Analog Devices"

# --- Row 8: HMC220B (re-uses the old placeholder formatted row) ----------
$ws.Range("A8").Value = "HMC220B"
$ws.Range("B8").Value = "8-MSOP-EP"
$ws.Range("C8").Value = "HMC220B 8-MSOP-EP"
$ws.Range("D8").Value = "HMC220B"
$ws.Range("E8").Value = "HMC220B"
$ws.Range("F8").Value = "8-MSOP-EP"
# G8 used to carry the old "Text" style with no content - drop that formatting
# before writing the manufacturer name so it ends up unstyled like the rest.
$ws.Range("G8").ClearFormats()
$ws.Range("G8").Value = "Analog Devices"
# I8 stays empty but picks up the new red-font highlight style.
$ws.Range("I8").Font.Color = 255

# --- Row 9: HMC369LP3 ------------------------------------------------------
$ws.Range("A9").Value = "HMC369LP3"
$ws.Range("B9").Value = "16-QFN (3X3)"
$ws.Range("C9").Value = "HMC369LP3 16-QFN (3X3)"
$ws.Range("D9").Value = "HMC369LP3"
$ws.Range("E9").Value = "HMC369LP3"
$ws.Range("F9").Value = "16-QFN (3X3)"
$ws.Range("G9").Value = "Analog Devices"

# --- Row 10: HMC451LP3 ------------------------------------------------------
$ws.Range("A10").Value = "HMC451LP3"
$ws.Range("B10").Value = "16-QFN (3X3)"
$ws.Range("C10").Value = "HMC451LP3 16-QFN (3X3)"
$ws.Range("D10").Value = "HMC451LP3"
$ws.Range("E10").Value = "HMC451LP3"
$ws.Range("F10").Value = "16-QFN (3X3)"
$ws.Range("G10").Value = "Analog Devices"

# --- Row 11: MTX2-73+ (Mini-Circuits) --------------------------------------
$ws.Range("A11").Value = "MTX2-73+"
$ws.Range("B11").Value = "12-QFN (3Х3)"
$ws.Range("C11").Value = "MTX2-73+ 12-QFN (3Х3)"
$ws.Range("D11").Value = "MTX2-73+"
$ws.Range("E11").Value = "MTX2-73+"
$ws.Range("F11").Value = "12-QFN (3Х3)"
$ws.Range("G11").Value = "Mimi-Circuits"

# --- Sheet-level bits -------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

$ws.Range("E8").Select()
